$d = $word.ActiveDocument

# --- Insert the new paragraph before the "Briefly describe the experiment" bullet ---
$anchor = $d.Content
$anchor.Find.Execute("• Briefly describe the experiment") | Out-Null
$anchor.Collapse(1)
$insertStart = $anchor.Start

$firstPart  = "Revising and extending p"
$secondPart = "revious drafts of nuclear and chloroblast genome sequences of Q. lobata "

# Insert the whole sentence (plus a trailing paragraph mark) as one shot, then
# re-anchor the "_GoBack" bookmark in the middle of it. Word treats "_GoBack"
# as a singleton bookmark, so adding it here automatically removes it from
# wherever it used to be (the end of the document).
$anchor.InsertBefore($firstPart + $secondPart + "`r")

$bmPos = $insertStart + $firstPart.Length
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
